$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Trial 1")
$ws2 = $wb.Worksheets.Item("Trial 2")

# --- Trial 2: negate the gauge-pressure column (B2:B43) ---
for ($r = 2; $r -le 43; $r++) {
    $cell = $ws2.Cells.Item($r, 2)
    $cell.Value = -1 * $cell.Value2()
}

# --- Trial 2: re-enter the C and D formulas as fill-down ranges so Excel
#     records them as shared formulas (matches the OOXML in the target) ---
$ws2.Range("C2:C43").Formula = "=A2+1"
$ws2.Range("D2:D43").Formula = "=B2+273.15"

# --- Selection / active-sheet bookkeeping ---
# Land on Trial 2 first and leave its last selection at O40 ...
$ws2.Activate() | Out-Null
$ws2.Range("O40").Select() | Out-Null

# ... then finish on Trial 1 with K38 selected, which is the sheet that
# ends up active/tabSelected when the file is saved.
$ws1.Activate() | Out-Null
$ws1.Range("K38").Select() | Out-Null
